# Journal de travail - Création du projet Unity / Implémentation du menu et de l'échiquier
#
# Fills in the journal entries for the Unity chess project creation, the
# temporary menu and the chessboard initialisation, then extends the
# "Tableau1" structured table so the new rows (and a few buffer rows) are
# included in it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rows 48:53 (and later 56) don't have any Date/Début/Fin cells yet, so pick
# up the existing "Date"/time number formatting from row 46 before filling
# in the values, so the new cells end up styled the same way as the rest of
# the table instead of using a generic default style.
$ws.Range("B46").Copy() | Out-Null
$ws.Range("B48:B53").PasteSpecial(-4122) | Out-Null
$ws.Range("C46:D46").Copy() | Out-Null
$ws.Range("C48:D53").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- New work-log rows (Tableau1 columns: Date, Début, Fin, Durée, Sujet, Description) ---

# Row 47 - 06.05.2024, Analyse
$ws.Range("B47").Value = 45418
$ws.Range("C47").Value = 0.57986111111111105
$ws.Range("D47").Value = 0.62847222222222221
$ws.Range("E47").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
$ws.Range("F47").Value = "Analyse"
$ws.Range("G47").Value = "Création des scénarios de tests"

# Row 48 - 07.05.2024, Analyse
$ws.Range("B48").Value = 45419
$ws.Range("C48").Value = 0.33333333333333331
$ws.Range("D48").Value = 0.3923611111111111
$ws.Range("E48").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
$ws.Range("F48").Value = "Analyse"
$ws.Range("G48").Value = "Réalisation des scénarios de tests"

# Rows 49 and 50 are filled in together, but the description of row 50 was
# actually typed first (this matches the original shared-string order).
$ws.Range("B49").Value = 45419
$ws.Range("C49").Value = 0.3923611111111111
$ws.Range("D49").Value = 0.39930555555555558
$ws.Range("E49").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
$ws.Range("F49").Value = "Implémentation"

$ws.Range("B50").Value = 45419
$ws.Range("C50").Value = 0.40972222222222227
$ws.Range("D50").Value = 0.46875
$ws.Range("E50").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
$ws.Range("F50").Value = "Implémentation"
$ws.Range("G50").Value = "Ajout d'un menu temporaire du jeu"

$ws.Range("G49").Value = "Ajout des scènes dans Unity"

# Row 51 - 07.05.2024, Implémentation
$ws.Range("B51").Value = 45419
$ws.Range("C51").Value = 0.46875
$ws.Range("D51").Value = 0.51041666666666663
$ws.Range("E51").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
$ws.Range("F51").Value = "Implémentation"
$ws.Range("G51").Value = "Création de la classe de l'échiquier et mise en place de la scène"

# Row 52 - 07.05.2024, Implémentation
$ws.Range("B52").Value = 45419
$ws.Range("C52").Value = 0.5625
$ws.Range("D52").Value = 0.61458333333333337
$ws.Range("E52").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
$ws.Range("F52").Value = "Implémentation"
$ws.Range("G52").Value = "Implémentation de l'initialisation de l'échiquier"

# Row 53 - 07.05.2024, Documentation (journal de travail)
$ws.Range("B53").Value = 45419
$ws.Range("C53").Value = 0.61458333333333337
$ws.Range("D53").Value = 0.62847222222222221
$ws.Range("E53").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
$ws.Range("F53").Value = "Documentation"

# --- Extend the structured table to keep some free rows below the data ---
$lo = $ws.ListObjects.Item("Tableau1")
$lo.Resize($ws.Range("B3:I63"))

# The calculated "Durée" column formula needs to be (re)applied to the
# newly added buffer rows of the table.
for ($r = 54; $r -le 63; $r++) {
    $ws.Range("E$r").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
}

# One date was already noted ahead of time on row 56
$ws.Range("B46").Copy() | Out-Null
$ws.Range("B56").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B56").Value = 45425

# --- Update the view so it matches where the user ended up working ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C55").Select()

$wb.Application.CalculateFull()
